$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value = 247
$ws1.Range("F3").Value = 1090
$ws1.Range("F5").Value = 430
$ws1.Range("F7").Value = 558
$ws1.Range("F8").Value = 70
$ws1.Range("F9").Value = 6804
$ws1.Range("F10").Value = 159
$ws1.Range("F12").Value = 143
$ws1.Range("F15").Value = 1099
$ws1.Range("F16").Value = 16203
$ws1.Range("F17").Value = 1590
$ws1.Range("F18").Value = 39
$ws1.Range("F20").Value = 186
$ws1.Range("F22").Value = 11358
$ws1.Range("F24").Value = 998
$ws1.Range("F25").Value = 4473
$ws1.Range("F26").Value = 316
$ws1.Range("F28").Value = 46
$ws1.Range("F29").Value = 845
$ws1.Range("F30").Value = 321

# Sheet 4: 全部类型
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F2").Value = 247
$ws4.Range("F3").Value = 1090
$ws4.Range("F5").Value = 430
$ws4.Range("F7").Value = 558
$ws4.Range("F9").Value = 70
$ws4.Range("F10").Value = 6804
$ws4.Range("F11").Value = 159
$ws4.Range("F13").Value = 143
$ws4.Range("F17").Value = 1099
$ws4.Range("F18").Value = 16203
$ws4.Range("F19").Value = 1590
$ws4.Range("F20").Value = 39
$ws4.Range("F22").Value = 186
$ws4.Range("F26").Value = 11358
$ws4.Range("F28").Value = 998
$ws4.Range("F29").Value = 4473
$ws4.Range("F30").Value = 316
$ws4.Range("F32").Value = 46
$ws4.Range("F33").Value = 845
$ws4.Range("F34").Value = 321
